$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking values
# ("5.09", "125.69", etc.) are preserved as strings, matching the
# source data which stores prices/volumes as text (e.g. "59.877.33").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.194.08"
$ws.Range("E2").Value = "  -4.63%  "
$ws.Range("D3").Value = "2.982.98"
$ws.Range("E3").Value = "  -6.08%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "579.78"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "125.69"
$ws.Range("E6").Value = "  -6.90%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "2.979.73"
$ws.Range("E8").Value = "  -6.12%  "
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("E10").Value = "  -5.86%  "
$ws.Range("D11").Value = "5.09"
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("D12").Value = "0.441"
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").Value = "0.0000222"
$ws.Range("E13").Value = "  -6.00%  "
$ws.Range("D14").Value = "32.46"
$ws.Range("E14").Value = "  -5.69%  "
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "3.466.37"
$ws.Range("E16").Value = "  -6.33%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "60.140.58"
$ws.Range("E17").Value = "  -4.58%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.970.46"
$ws.Range("E18").Value = "  -6.47%  "
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  -5.27%  "
$ws.Range("D20").Value = "433.15"
$ws.Range("E20").Value = "  -6.14%  "
$ws.Range("D21").Value = "13.09"
$ws.Range("E21").Value = "  -6.38%  "
$ws.Range("D22").Value = "0.662"
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("D23").Value = "7.01"
$ws.Range("E23").Value = "  -7.72%  "
$ws.Range("D24").Value = "12.62"
$ws.Range("E24").Value = "  -5.09%  "
$ws.Range("D25").Value = "79.09"
$ws.Range("E25").Value = "  -3.95%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "2.54"
$ws.Range("E28").Value = "  -4.81%  "
$ws.Range("D29").Value = "7.25"
$ws.Range("E29").Value = "  -5.10%  "
$ws.Range("D30").Value = "1.88"
$ws.Range("E30").Value = "  -7.21%  "
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  -8.66%  "
$ws.Range("D32").Value = "25.30"
$ws.Range("E32").Value = "  -6.87%  "
$ws.Range("D33").Value = "0.0937"
$ws.Range("E33").Value = "  -7.79%  "
$ws.Range("D34").Value = "2.18"
$ws.Range("E34").Value = "  -7.77%  "
$ws.Range("D35").Value = "0.948"
$ws.Range("E35").Value = "  -8.08%  "
$ws.Range("D36").Value = "5.59"
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("D37").Value = "49.42"
$ws.Range("E37").Value = "  -3.72%  "
$ws.Range("D38").Value = "0.0₃0660"
$ws.Range("E38").Value = "  -6.52%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0360"
$ws.Range("E39").Value = "  -7.10%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "7.94"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "384.76"
$ws.Range("E41").Value = "  -4.69%  "
$ws.Range("D42").Value = "0.109"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").Value = "2.46"
$ws.Range("E43").Value = "  -6.86%  "
$ws.Range("D44").Value = "2.629.79"
$ws.Range("E44").Value = "  -6.52%  "
$ws.Range("D46").Value = "0.236"
$ws.Range("E46").Value = "  -6.27%  "
$ws.Range("D47").Value = "119.05"
$ws.Range("E47").Value = "  -4.08%  "
$ws.Range("D48").Value = "1.98"
$ws.Range("E48").Value = "  -5.73%  "
$ws.Range("E49").Value = "  -3.69%  "
$ws.Range("D50").Value = "23.40"
$ws.Range("E50").Value = "  -6.80%  "
$ws.Range("D51").Value = "31.33"
$ws.Range("E51").Value = "  -9.68%  "
